# Insert a new weekly price record for "Macroferia Regional de Talca - Repollo".
# The new record belongs right after the existing row 504 (chronologically it
# sits between the current rows 504 and 505), so we insert a brand-new row at
# position 505, which pushes the former rows 505:524 down to 506:525.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 505 (shifts 505:524 -> 506:525).
$ws.Rows.Item(505).Insert()

# Populate the new row with the new record's data.
$ws.Range("A505").Value = 5
$ws.Range("B505").Value = "Macroferia Regional de Talca"
$ws.Range("C505").Value = "Maule"
$ws.Range("D505").Value = 45147
$ws.Range("E505").Value = 7
$ws.Range("F505").Value = 100112006
$ws.Range("G505").Value = "Repollo"
$ws.Range("H505").Value = "Crespo record"
$ws.Range("I505").Value = "Primera"
$ws.Range("J505").Value = 5000
$ws.Range("K505").Value = 500
$ws.Range("L505").Value = 500
$ws.Range("M505").Value = 500
$ws.Range("N505").Value = "`$/unidad"
$ws.Range("O505").Value = "Región del Maule"
$ws.Range("P505").Value = 500
$ws.Range("Q505").Value = 1
$ws.Range("R505").Value = "Hortaliza"
